# Fix import controller: the header labels in the VUNG lookup sheet used
# internal field names (MaVung / TenVung) instead of the display captions
# the DB import expects (MÃ VÙNG / TÊN VÙNG), which caused saves to fail.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VUNG")

$ws.Range("A1").Value = "MÃ VÙNG"
$ws.Range("B1").Value = "TÊN VÙNG"
